$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the year value next to "Pune" (A1), stored as text "2025" (leading
# apostrophe forces text, matching the target's string-typed cell).
$ws.Range("B1").Value = "'2025"
